# Restaurant Manager.pptx -- "Updates to class diagram."
#
# Slide 4 ("Requirements") content placeholder loses the
# "Home page is to have a message area" bullet, and (since the
# remaining text now fits the placeholder) the body's autofit no
# longer needs to shrink the font -- it goes from
# <a:normAutofit fontScale="92500" lnSpcReduction="10000"/> back to a
# plain <a:normAutofit/>.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(4)

$contentShape = $null
foreach ($shp in $slide.Shapes) {
    if ($shp.Name -eq "Content Placeholder 2") {
        $contentShape = $shp
    }
}

$tr = $contentShape.TextFrame.TextRange
$bullets = @(
    "Home page is to have a summary of locations assigned to employee",
    "Show current location and employee info",
    "Home page is to show business analytics",
    "Home page is to show employee schedule",
    "A page for parties",
    "Placement of customer orders"
)
$tr.Text = [string]::Join([char]13, $bullets)

# Re-apply "shrink text on overflow" now that the bullet list is
# shorter; it now fits without scaling, so PowerPoint drops the
# fontScale/lnSpcReduction attributes.
$contentShape.TextFrame.AutoSize = 2
